$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): force text format so numeric-looking strings
# such as "0.171" or "69.174.53" are preserved exactly as text,
# matching the original inline-string cell type.
$dCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D17","D21","D22","D23","D24","D25","D27","D31","D32","D33","D34","D35","D38","D39","D41","D42","D43","D44","D46","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.174.53"
$ws.Range("D3").Value = "3.752.30"
$ws.Range("D5").Value = "602.22"
$ws.Range("D6").Value = "167.11"
$ws.Range("D7").Value = "3.750.22"
$ws.Range("D9").Value = "0.539"
$ws.Range("D10").Value = "0.170"
$ws.Range("D11").Value = "6.39"
$ws.Range("D13").Value = "37.97"
$ws.Range("D14").Value = "0.0000248"
$ws.Range("D15").Value = "4.378.85"
$ws.Range("D16").Value = "3.756.55"
$ws.Range("D17").Value = "69.168.92"
$ws.Range("D21").Value = "11.10"
$ws.Range("D22").Value = "493.25"
$ws.Range("D23").Value = "0.728"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("D25").Value = "84.91"
$ws.Range("D27").Value = "12.29"
$ws.Range("D31").Value = "2.47"
$ws.Range("D32").Value = "8.12"
$ws.Range("D33").Value = "31.53"
$ws.Range("D34").Value = "3.899.28"
$ws.Range("D35").Value = "3.684.93"
$ws.Range("D38").Value = "5.98"
$ws.Range("D39").Value = "1.01"
$ws.Range("D41").Value = "0.326"
$ws.Range("D42").Value = "3.03"
$ws.Range("D43").Value = "48.67"
$ws.Range("D44").Value = "426.03"
$ws.Range("D46").Value = "8.47"
$ws.Range("D48").Value = "40.18"
$ws.Range("D49").Value = "141.50"
$ws.Range("D50").Value = "2.793.28"

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Other columns (B, C, E): plain text values, no numeric coercion risk
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("E21").Value = "  +8.30%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("E24").Value = "  +7.58%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +5.70%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("E51").Value = "  +0.59%  "
